$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ A=1.8; B=0.3; C=200; D=25; E=0; F=8;  G="nach ca 40 M Comp. => 1.0 rating" },
    @{ A=1.8; B=0.3; C=150; D=25; E=0; F=9;  G="nach ca 80 M Comp. => 0.99 rating" },
    @{ A=1.8; B=0.3; C=100; D=25; E=0; F=$null; G="nach ca 10 M Comp. => stagnation 0.5 rating" },
    @{ A=1.8; B=0.3; C=300; D=25; E=0; F=10; G="nach ca 38 M Comp. => 1.0 rating" },
    @{ A=1.8; B=0.3; C=350; D=25; E=0; F=11; G="nach ca 60 M Comp. => 1.0 rating" }
)

$r = 12
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 1).NumberFormat = "0.00"
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    if ($row.F -ne $null) {
        $ws.Cells.Item($r, 6).Value = $row.F
    }
    $ws.Cells.Item($r, 7).Value = $row.G
    $r++
}

$ws.Range("F16").Select()
